$wb = $excel.ActiveWorkbook

# --- Tournament: drop the competition-key / host-key rows (2:3) and the ---
# --- venue-key.N rows (now 15:24 after the first delete) ------------------
$tournament = $wb.Worksheets.Item("Tournament")
$tournament.Range("A2:A3").EntireRow.Delete()
$tournament.Range("A15:A24").EntireRow.Select()
$tournament.Range("A15:A24").EntireRow.Delete()

# --- Colors: just a leftover selection change in the saved view ----------
$colors = $wb.Worksheets.Item("Colors")
$colors.Activate()
$colors.Range("K1").EntireColumn.Select()

# --- new Properties sheet, appended after #venues -------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$props = $wb.Worksheets.Add($null, $lastSheet)
$props.Name = "Properties"

# final row -> content, keyed by the row each entry lands on once the table
# is sorted A-Z by key (competition/host/timezone/color.*/venue.*)
$finalRows = @{
   1 = @("key",         "value",           "notes")
   2 = @("competition", "mens-euro",       "")
   3 = @("host",        "germany",         "")
   4 = @("timezone",    "Europe/Berlin",   "")
   5 = @("color.a",     "#c4e1b5",         "pale green")
   6 = @("color.b",     "#b0d0ee",         "pale blue")
   7 = @("color.c",     "#f79d8f",         "pale red")
   8 = @("color.d",     "#fee289",         "pale yellow")
   9 = @("color.e",     "#c0e4df",         "pale teal")
  10 = @("color.f",     "#acacac",         "pale gray")
  11 = @("venue.01",    "de-berlin",       "")
  12 = @("venue.02",    "de-cologne",      "")
  13 = @("venue.03",    "de-dortmund",     "")
  14 = @("venue.04",    "de-dusseldorf",   "")
  15 = @("venue.05",    "de-frankfurt",    "")
  16 = @("venue.06",    "de-gelsenkirchen","")
  17 = @("venue.07",    "de-hamburg",      "")
  18 = @("venue.08",    "de-leipzig",      "")
  19 = @("venue.09",    "de-munich",       "")
  20 = @("venue.10",    "de-stuttgart",    "")
}

# write in the same order the strings were originally authored (header,
# then the colors, then the venues, then competition/host, then timezone)
# so the shared-string table gets the new entries in that order
$writeOrder = @(1, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 2, 3, 4)

foreach ($r in $writeOrder) {
  $row = $finalRows[$r]
  for ($c = 0; $c -lt $row.Count; $c++) {
    $value = $row[$c]
    if ($value -ne "") {
      $props.Cells.Item($r, $c + 1).Value = $value
    }
  }
}

$tbl = $props.ListObjects.Add(1, $props.Range("A1:C20"), $null, 1)
$tbl.Name = "Properties"

$props.Columns.Item(1).EntireColumn.AutoFit()
$props.Columns.Item(2).EntireColumn.AutoFit()
$props.Columns.Item(3).EntireColumn.AutoFit()

$props.Range("A4:B4").Select()
$props.Activate()
